$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Date) values look like dates ("2025-12-17") and Excel would
# auto-convert them to date serials on assignment. Force text interpretation
# by briefly setting a text number format, then clear the format again so the
# cell keeps its default (unstyled) look, matching the source data which is
# plain text, not a real date value.

# Row 2
$ws.Cells.Item(2, 1).Value = 'Turkish 2 Lig'
$ws.Cells.Item(2, 2).NumberFormat = '@'
$ws.Cells.Item(2, 2).Value = '2025-12-17'
$ws.Cells.Item(2, 2).ClearFormats()
$ws.Cells.Item(2, 3).Value = '08:00:00'
$ws.Cells.Item(2, 4).Value = 'Iskenderunspor'
$ws.Cells.Item(2, 5).Value = 'Adana 1954 FK'
$ws.Cells.Item(2, 6).Value = 1.02
$ws.Cells.Item(2, 7).Value = 1000
$ws.Cells.Item(2, 8).Value = 1.02
$ws.Cells.Item(2, 9).Value = 1000
$ws.Cells.Item(2, 10).Value = 1.02
$ws.Cells.Item(2, 11).Value = 1000
$ws.Cells.Item(2, 12).Value = 1.01
$ws.Cells.Item(2, 13).Value = 1.01
$ws.Cells.Item(2, 14).Value = 1.25
$ws.Cells.Item(2, 15).Value = 1.01
$ws.Cells.Item(2, 16).Value = 1.24
$ws.Cells.Item(2, 17).Value = 1.02
$ws.Cells.Item(2, 18).Value = 1.18
$ws.Cells.Item(2, 19).Value = 1.36
$ws.Cells.Item(2, 20).Value = 1.01
$ws.Cells.Item(2, 21).Value = 1.01
$ws.Cells.Item(2, 22).Value = 1.01
$ws.Cells.Item(2, 23).Value = 1.01
$ws.Cells.Item(2, 24).Value = 1000
$ws.Cells.Item(2, 25).Value = 1000
$ws.Cells.Item(2, 26).Value = 1000
$ws.Cells.Item(2, 27).Value = 1000
$ws.Cells.Item(2, 28).Value = 1000
$ws.Cells.Item(2, 29).Value = 1000
$ws.Cells.Item(2, 30).Value = 1000
$ws.Cells.Item(2, 31).Value = 1000
$ws.Cells.Item(2, 32).Value = 1000
$ws.Cells.Item(2, 33).Value = 1000
$ws.Cells.Item(2, 34).Value = 1000
$ws.Cells.Item(2, 35).Value = 1000
$ws.Cells.Item(2, 36).Value = 1000
$ws.Cells.Item(2, 37).Value = 1000
$ws.Cells.Item(2, 38).Value = 1000
$ws.Cells.Item(2, 39).Value = 1000
$ws.Cells.Item(2, 40).Value = 1000
$ws.Cells.Item(2, 41).Value = 1000

# Row 3
$ws.Cells.Item(3, 1).Value = 'Turkish 2 Lig'
$ws.Cells.Item(3, 2).NumberFormat = '@'
$ws.Cells.Item(3, 2).Value = '2025-12-17'
$ws.Cells.Item(3, 2).ClearFormats()
$ws.Cells.Item(3, 3).Value = '09:00:00'
$ws.Cells.Item(3, 4).Value = 'Beykoz Anadolu Spor'
$ws.Cells.Item(3, 5).Value = 'Batman Petrolspor'
$ws.Cells.Item(3, 6).Value = 1.02
$ws.Cells.Item(3, 7).Value = 1000
$ws.Cells.Item(3, 8).Value = 1.02
$ws.Cells.Item(3, 9).Value = 1000
$ws.Cells.Item(3, 10).Value = 1.02
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 12).Value = 1.01
$ws.Cells.Item(3, 13).Value = 1.01
$ws.Cells.Item(3, 14).Value = 1.24
$ws.Cells.Item(3, 15).Value = 1.01
$ws.Cells.Item(3, 16).Value = 1.24
$ws.Cells.Item(3, 17).Value = 1.01
$ws.Cells.Item(3, 18).Value = 1.18
$ws.Cells.Item(3, 19).Value = 1.02
$ws.Cells.Item(3, 20).Value = 1.01
$ws.Cells.Item(3, 21).Value = 1.01
$ws.Cells.Item(3, 22).Value = 1.01
$ws.Cells.Item(3, 23).Value = 1.01
$ws.Cells.Item(3, 24).Value = 1000
$ws.Cells.Item(3, 25).Value = 1000
$ws.Cells.Item(3, 26).Value = 1000
$ws.Cells.Item(3, 27).Value = 1000
$ws.Cells.Item(3, 28).Value = 1000
$ws.Cells.Item(3, 29).Value = 1000
$ws.Cells.Item(3, 30).Value = 1000
$ws.Cells.Item(3, 31).Value = 1000
$ws.Cells.Item(3, 32).Value = 1000
$ws.Cells.Item(3, 33).Value = 1000
$ws.Cells.Item(3, 34).Value = 1000
$ws.Cells.Item(3, 35).Value = 1000
$ws.Cells.Item(3, 36).Value = 1000
$ws.Cells.Item(3, 37).Value = 1000
$ws.Cells.Item(3, 38).Value = 1000
$ws.Cells.Item(3, 39).Value = 1000
$ws.Cells.Item(3, 40).Value = 1000
$ws.Cells.Item(3, 41).Value = 1000

# Row 4
$ws.Cells.Item(4, 1).Value = 'Serbian Super League'
$ws.Cells.Item(4, 2).NumberFormat = '@'
$ws.Cells.Item(4, 2).Value = '2025-12-17'
$ws.Cells.Item(4, 2).ClearFormats()
$ws.Cells.Item(4, 3).Value = '12:30:00'
$ws.Cells.Item(4, 4).Value = 'FK Radnicki 1923'
$ws.Cells.Item(4, 5).Value = 'Cukaricki'
$ws.Cells.Item(4, 6).Value = 1.04
$ws.Cells.Item(4, 7).Value = 1000
$ws.Cells.Item(4, 8).Value = 1.28
$ws.Cells.Item(4, 9).Value = 1000
$ws.Cells.Item(4, 10).Value = 1.02
$ws.Cells.Item(4, 11).Value = 4.7
$ws.Cells.Item(4, 12).Value = 1.01
$ws.Cells.Item(4, 13).Value = 1.04
$ws.Cells.Item(4, 14).Value = 1.21
$ws.Cells.Item(4, 15).Value = 1.21
$ws.Cells.Item(4, 16).Value = 1.21
$ws.Cells.Item(4, 17).Value = 1.21
$ws.Cells.Item(4, 18).Value = 1.09
$ws.Cells.Item(4, 19).Value = 1.01
$ws.Cells.Item(4, 20).Value = 1.01
$ws.Cells.Item(4, 21).Value = 1.01
$ws.Cells.Item(4, 22).Value = 1.01
$ws.Cells.Item(4, 23).Value = 1.01
$ws.Cells.Item(4, 24).Value = 1000
$ws.Cells.Item(4, 25).Value = 1000
$ws.Cells.Item(4, 26).Value = 1000
$ws.Cells.Item(4, 27).Value = 1000
$ws.Cells.Item(4, 28).Value = 1000
$ws.Cells.Item(4, 29).Value = 1000
$ws.Cells.Item(4, 30).Value = 1000
$ws.Cells.Item(4, 31).Value = 1000
$ws.Cells.Item(4, 32).Value = 1000
$ws.Cells.Item(4, 33).Value = 1000
$ws.Cells.Item(4, 34).Value = 1000
$ws.Cells.Item(4, 35).Value = 1000
$ws.Cells.Item(4, 36).Value = 1000
$ws.Cells.Item(4, 37).Value = 1000
$ws.Cells.Item(4, 38).Value = 1000
$ws.Cells.Item(4, 39).Value = 1000
$ws.Cells.Item(4, 40).Value = 1000
$ws.Cells.Item(4, 41).Value = 1000

# Row 5
$ws.Cells.Item(5, 1).Value = 'Swiss Super League'
$ws.Cells.Item(5, 2).NumberFormat = '@'
$ws.Cells.Item(5, 2).Value = '2025-12-17'
$ws.Cells.Item(5, 2).ClearFormats()
$ws.Cells.Item(5, 3).Value = '16:30:00'
$ws.Cells.Item(5, 4).Value = 'Young Boys'
$ws.Cells.Item(5, 5).Value = 'Grasshoppers Zurich'
$ws.Cells.Item(5, 6).Value = 1.6
$ws.Cells.Item(5, 7).Value = 1.63
$ws.Cells.Item(5, 8).Value = 5.3
$ws.Cells.Item(5, 9).Value = 5.8
$ws.Cells.Item(5, 10).Value = 5
$ws.Cells.Item(5, 11).Value = 5.1
$ws.Cells.Item(5, 12).Value = 1.26
$ws.Cells.Item(5, 13).Value = 1.03
$ws.Cells.Item(5, 14).Value = 5
$ws.Cells.Item(5, 15).Value = 1.17
$ws.Cells.Item(5, 16).Value = 2.46
$ws.Cells.Item(5, 17).Value = 1.48
$ws.Cells.Item(5, 18).Value = 1.72
$ws.Cells.Item(5, 19).Value = 1.98
$ws.Cells.Item(5, 20).Value = 1.6
$ws.Cells.Item(5, 21).Value = 2.46
$ws.Cells.Item(5, 22).Value = 1.21
$ws.Cells.Item(5, 23).Value = 2.58
$ws.Cells.Item(5, 24).Value = 29
$ws.Cells.Item(5, 25).Value = 34
$ws.Cells.Item(5, 26).Value = 55
$ws.Cells.Item(5, 27).Value = 160
$ws.Cells.Item(5, 28).Value = 15
$ws.Cells.Item(5, 29).Value = 12
$ws.Cells.Item(5, 30).Value = 980
$ws.Cells.Item(5, 31).Value = 60
$ws.Cells.Item(5, 32).Value = 13
$ws.Cells.Item(5, 33).Value = 10.5
$ws.Cells.Item(5, 34).Value = 17.5
$ws.Cells.Item(5, 35).Value = 55
$ws.Cells.Item(5, 36).Value = 17.5
$ws.Cells.Item(5, 37).Value = 15
$ws.Cells.Item(5, 38).Value = 25
$ws.Cells.Item(5, 39).Value = 70
$ws.Cells.Item(5, 40).Value = 6.8
$ws.Cells.Item(5, 41).Value = 44

# Row 6
$ws.Cells.Item(6, 1).Value = 'Swiss Super League'
$ws.Cells.Item(6, 2).NumberFormat = '@'
$ws.Cells.Item(6, 2).Value = '2025-12-17'
$ws.Cells.Item(6, 2).ClearFormats()
$ws.Cells.Item(6, 3).Value = '16:30:00'
$ws.Cells.Item(6, 4).Value = 'FC Zurich'
$ws.Cells.Item(6, 5).Value = 'Lugano'
$ws.Cells.Item(6, 6).Value = 2.76
$ws.Cells.Item(6, 7).Value = 3
$ws.Cells.Item(6, 8).Value = 2.46
$ws.Cells.Item(6, 9).Value = 2.72
$ws.Cells.Item(6, 10).Value = 3.55
$ws.Cells.Item(6, 11).Value = 4
$ws.Cells.Item(6, 12).Value = 1.29
$ws.Cells.Item(6, 13).Value = 1.05
$ws.Cells.Item(6, 14).Value = 4.4
$ws.Cells.Item(6, 15).Value = 1.24
$ws.Cells.Item(6, 16).Value = 2.16
$ws.Cells.Item(6, 17).Value = 1.63
$ws.Cells.Item(6, 18).Value = 1.46
$ws.Cells.Item(6, 19).Value = 2.8
$ws.Cells.Item(6, 20).Value = 1.64
$ws.Cells.Item(6, 21).Value = 2.32
$ws.Cells.Item(6, 22).Value = 1.58
$ws.Cells.Item(6, 23).Value = 1.5
$ws.Cells.Item(6, 24).Value = 22
$ws.Cells.Item(6, 25).Value = 15.5
$ws.Cells.Item(6, 26).Value = 23
$ws.Cells.Item(6, 27).Value = 44
$ws.Cells.Item(6, 28).Value = 17
$ws.Cells.Item(6, 29).Value = 10.5
$ws.Cells.Item(6, 30).Value = 15
$ws.Cells.Item(6, 31).Value = 32
$ws.Cells.Item(6, 32).Value = 26
$ws.Cells.Item(6, 33).Value = 16
$ws.Cells.Item(6, 34).Value = 19.5
$ws.Cells.Item(6, 35).Value = 42
$ws.Cells.Item(6, 36).Value = 55
$ws.Cells.Item(6, 37).Value = 38
$ws.Cells.Item(6, 38).Value = 980
$ws.Cells.Item(6, 39).Value = 85
$ws.Cells.Item(6, 40).Value = 27
$ws.Cells.Item(6, 41).Value = 22

# Row 7
$ws.Cells.Item(7, 1).Value = 'Swiss Super League'
$ws.Cells.Item(7, 2).NumberFormat = '@'
$ws.Cells.Item(7, 2).Value = '2025-12-17'
$ws.Cells.Item(7, 2).ClearFormats()
$ws.Cells.Item(7, 3).Value = '16:30:00'
$ws.Cells.Item(7, 4).Value = 'Luzern'
$ws.Cells.Item(7, 5).Value = 'FC Basel'
$ws.Cells.Item(7, 6).Value = 3.5
$ws.Cells.Item(7, 7).Value = 3.75
$ws.Cells.Item(7, 8).Value = 2.08
$ws.Cells.Item(7, 9).Value = 2.14
$ws.Cells.Item(7, 10).Value = 3.85
$ws.Cells.Item(7, 11).Value = 4.1
$ws.Cells.Item(7, 12).Value = 1.27
$ws.Cells.Item(7, 13).Value = 1.03
$ws.Cells.Item(7, 14).Value = 5.7
$ws.Cells.Item(7, 15).Value = 1.18
$ws.Cells.Item(7, 16).Value = 2.58
$ws.Cells.Item(7, 17).Value = 1.55
$ws.Cells.Item(7, 18).Value = 1.65
$ws.Cells.Item(7, 19).Value = 2.32
$ws.Cells.Item(7, 20).Value = 1.53
$ws.Cells.Item(7, 21).Value = 2.64
$ws.Cells.Item(7, 22).Value = 1.88
$ws.Cells.Item(7, 23).Value = 1.37
$ws.Cells.Item(7, 24).Value = 1000
$ws.Cells.Item(7, 25).Value = 18.5
$ws.Cells.Item(7, 26).Value = 1000
$ws.Cells.Item(7, 27).Value = 980
$ws.Cells.Item(7, 28).Value = 1000
$ws.Cells.Item(7, 29).Value = 10
$ws.Cells.Item(7, 30).Value = 12
$ws.Cells.Item(7, 31).Value = 1000
$ws.Cells.Item(7, 32).Value = 1000
$ws.Cells.Item(7, 33).Value = 16
$ws.Cells.Item(7, 34).Value = 18
$ws.Cells.Item(7, 35).Value = 1000
$ws.Cells.Item(7, 36).Value = 60
$ws.Cells.Item(7, 37).Value = 980
$ws.Cells.Item(7, 38).Value = 1000
$ws.Cells.Item(7, 39).Value = 55
$ws.Cells.Item(7, 40).Value = 1000
$ws.Cells.Item(7, 41).Value = 9.6

# Row 8
$ws.Cells.Item(8, 1).Value = 'Scottish Premiership'
$ws.Cells.Item(8, 2).NumberFormat = '@'
$ws.Cells.Item(8, 2).Value = '2025-12-17'
$ws.Cells.Item(8, 2).ClearFormats()
$ws.Cells.Item(8, 3).Value = '17:00:00'
$ws.Cells.Item(8, 4).Value = 'Dundee Utd'
$ws.Cells.Item(8, 5).Value = 'Celtic'
$ws.Cells.Item(8, 6).Value = 7.2
$ws.Cells.Item(8, 7).Value = 8
$ws.Cells.Item(8, 8).Value = 1.5
$ws.Cells.Item(8, 9).Value = 1.55
$ws.Cells.Item(8, 10).Value = 4.6
$ws.Cells.Item(8, 11).Value = 5.2
$ws.Cells.Item(8, 12).Value = 1.3
$ws.Cells.Item(8, 13).Value = 1.03
$ws.Cells.Item(8, 14).Value = 5.1
$ws.Cells.Item(8, 15).Value = 1.2
$ws.Cells.Item(8, 16).Value = 2.42
$ws.Cells.Item(8, 17).Value = 1.61
$ws.Cells.Item(8, 18).Value = 1.57
$ws.Cells.Item(8, 19).Value = 2.5
$ws.Cells.Item(8, 20).Value = 1.8
$ws.Cells.Item(8, 21).Value = 2.1
$ws.Cells.Item(8, 22).Value = 2.8
$ws.Cells.Item(8, 23).Value = 1.14
$ws.Cells.Item(8, 24).Value = 23
$ws.Cells.Item(8, 25).Value = 11.5
$ws.Cells.Item(8, 26).Value = 10.5
$ws.Cells.Item(8, 27).Value = 14.5
$ws.Cells.Item(8, 28).Value = 32
$ws.Cells.Item(8, 29).Value = 11.5
$ws.Cells.Item(8, 30).Value = 13
$ws.Cells.Item(8, 31).Value = 14.5
$ws.Cells.Item(8, 32).Value = 60
$ws.Cells.Item(8, 33).Value = 27
$ws.Cells.Item(8, 34).Value = 22
$ws.Cells.Item(8, 35).Value = 38
$ws.Cells.Item(8, 36).Value = 220
$ws.Cells.Item(8, 37).Value = 95
$ws.Cells.Item(8, 38).Value = 85
$ws.Cells.Item(8, 39).Value = 110
$ws.Cells.Item(8, 40).Value = 100
$ws.Cells.Item(8, 41).Value = 7
